# Lab 5 Task 0 finished
#
# Highlights the "7) Instruction memory" and "8) Data memory" sections
# (heading + body paragraphs) in yellow, the way Word does when the
# paragraph text is selected and a highlight color is applied (which
# also coalesces the runs that make up the paragraph into fewer, larger
# runs sharing the same formatting).
#
# A couple of the phrases being searched for recur several times in the
# document (e.g. "Post-synthesis Functional simulation", "the given .v
# file in the folder" ...), so every Find is anchored to start at (or
# after) the end of the previous match instead of always re-scanning
# from the top of the document - that way the right occurrence is
# always the one that gets touched.

$d = $word.ActiveDocument
$wdYellow = 7
$q1 = [char]0x201C   # “
$q2 = [char]0x201D   # ”
$docEnd = $d.Content.End

# Search for $text starting at character position $from (to the end of
# the document). Returns the matched Range (or $null if not found).
function Find-From([int]$from, [string]$text) {
    $r = $d.Range($from, $docEnd)
    $found = $r.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        return $null
    }
    return $r
}

# Replaces $text (found starting at $from) with itself, which merges
# the runs spanned by the match into a single run - same effect Word's
# own Find & Replace has. Returns the end position of the match.
function Merge-From([int]$from, [string]$text) {
    $r = $d.Range($from, $docEnd)
    [void]$r.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, $text, 2)
    return $from + $text.Length
}

# Anchor just before the section we are editing ("7) (10 pts)
# Instruction memory" heading) so every subsequent search is scoped to
# start there or later.
$anchor = Find-From 0 "then ALUResult = 32  (there are 32 0"
$pos = $anchor.End

# ---------------------------------------------------------------
# 7) (10 pts) Instruction memory  -- heading paragraph
# ---------------------------------------------------------------
$pos = Merge-From $pos "7) (10 pts) "
$rng = Find-From $anchor.End "7) (10 pts) Instruction memory"
$rng.HighlightColorIndex = $wdYellow
$pos = $rng.End

# ---------------------------------------------------------------
# "It is used to keep the machine language ..." -- body paragraph
# ---------------------------------------------------------------
$p109 = "It is used to keep the machine language (binary sequence of instructions that we want the processor to execute). The " + $q1 + "initial" + $q2 + " part of the code already contains the code that you will use in the next task of lab 5."
$start109 = $pos
$pos = Merge-From $pos $p109
$rng = Find-From $start109 $p109
$rng.HighlightColorIndex = $wdYellow
$pos = $rng.End

# ---------------------------------------------------------------
# "First, read the given .v file ... Post-synthesis Functional
# simulation on your component." -- body paragraph with an
# underlined phrase in the middle.
# ---------------------------------------------------------------
$p110pre = "First, read the given .v file in the folder " + $q1 + "InstructionMemory" + $q2 + " to understand how it should be used/how it functions. Then complete the provided testbench (use the waveform shown below for inputs) and "
$start110 = $pos
$pos = Merge-From $pos $p110pre

$p110tail = " on your component."
$rngAll = Find-From $start110 "run Post-synthesis Functional simulation on your component."
$tailStart = $rngAll.End - $p110tail.Length
$tailRng = $d.Range($tailStart, $rngAll.End)
$tailRng.Text = $p110tail

$p110full = $p110pre + "run Post-synthesis Functional simulation" + $p110tail
$rng = Find-From $start110 $p110full
$rng.HighlightColorIndex = $wdYellow
$pos = $rng.End

# ---------------------------------------------------------------
# 8) (10 pts) Data memory  -- heading paragraph
# ---------------------------------------------------------------
$start8 = $pos
$pos = Merge-From $pos "8) (10 pts) "
$pos = Merge-From $pos "Data memory"
$rng = Find-From $start8 "8) (10 pts) Data memory"
$rng.HighlightColorIndex = $wdYellow
$pos = $rng.End

# ---------------------------------------------------------------
# "It is used to keep the result from ALU ... Post-synthesis
# Functional simulation on your component." -- body paragraph. The
# original "Then complete ... and run " phrase ends up tucked away
# inside a (hidden) field in the target document, so it disappears
# from the visible/extracted text; "Post-synthesis Functional
# simulation" remains as a normal, underlined run.
# ---------------------------------------------------------------
$p114pre = "It is used to keep the result from ALU (used in the later tasks of lab 5). First, read the given .v file in the folder " + $q1 + "DataMemory" + $q2 + " to understand how it should be used/how it functions. "
$p114mid = "Then complete the provided testbench (use the waveform shown below for inputs) and run "
$p114us  = "Post-synthesis Functional simulation"
$p114tail = " on your component."

$start114 = $pos
$p114old = $p114pre + $p114mid + $p114us + $p114tail
$rngOld = Find-From $start114 $p114old
$rngOld.Text = $p114pre + $p114us + $p114tail

$rng = Find-From $start114 $p114pre
$rng.HighlightColorIndex = $wdYellow
$posUs = $rng.End

$rng = Find-From $posUs $p114us
$rng.Underline = 1
$rng.HighlightColorIndex = $wdYellow
$posTail = $rng.End

$rng = Find-From $posTail $p114tail
$rng.HighlightColorIndex = $wdYellow
